# Sprint 5 - "Casos de prueba": Daily 2 update.
# Changes the "Estado" (status) column of several test-case rows in the
# single table of the document from "Pendiente" to "Aprobado".
# Row FOR-943.1 additionally gets a detailed split: "Aprobado en backend"
# stays as the first paragraph and a new second paragraph
# "Pendiente detalle en front" is appended in the same cell.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Locate a data row by the "ID" shown in its first column (e.g. "FOR-943.1").
function Get-RowByCaseId($table, $caseId) {
    for ($r = 1; $r -le $table.Rows.Count; $r++) {
        if ($table.Cell($r, 1).Range.Text -like ("*" + $caseId + "*")) {
            return $r
        }
    }
    return -1
}

# --- Row FOR-943.1: split the status cell into two paragraphs ---
$row943 = Get-RowByCaseId $tbl "FOR-943.1"
$cell = $tbl.Cell($row943, 5)
$rng = $cell.Range

$newCellXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
  '<w:body>' + `
  '<w:p w14:paraId="06664B1A" w14:textId="77777777" w:rsidR="0042685F" w:rsidRPr="0042685F" w:rsidRDefault="0042685F" w:rsidP="00D965CC">' + `
    '<w:pPr><w:spacing w:before="240" w:after="240"/><w:rPr><w:bCs/></w:rPr></w:pPr>' + `
    '<w:r w:rsidRPr="0042685F"><w:rPr><w:bCs/></w:rPr><w:t>Aprobado</w:t></w:r>' + `
    '<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> en </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:bCs/></w:rPr><w:t>backend</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>' + `
  '<w:p>' + `
    '<w:pPr><w:spacing w:before="240" w:after="240"/><w:rPr><w:bCs/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">Pendiente detalle en </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:bCs/></w:rPr><w:t>front</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

[void]$rng.InsertXML($newCellXml)

# --- Simple "Pendiente" -> "Aprobado" rows ---
$caseIdsToApprove = @("FOR-944.1", "FOR-959.1", "FOR-973.1", "FOR-975.1", "FOR-977.1")

foreach ($caseId in $caseIdsToApprove) {
    $r = Get-RowByCaseId $tbl $caseId
    if ($r -gt 0) {
        $statusRange = $tbl.Cell($r, 5).Range
        $found = $statusRange.Find.Execute("Pendiente", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $statusRange.Text = "Aprobado"
        }
    }
}

Write-Output "done"
